$d = $word.ActiveDocument

# For each of the 10 "MATERIALS REQUIRED BUT NOT PROVIDED" bullet placeholders,
# replace the single Jinja placeholder run with two bullet-character runs
# (the first run carries an explicit 22-half-point font size), and set the
# paragraph's pPr to use explicit left justification. The indentation
# (w:ind left=360 firstLine=0) is applied in a second pass below because it
# must be set via ParagraphFormat after the paragraph text/runs are replaced.

for ($i = 1; $i -le 10; $i++) {
    $needle = "{{ req_material_$i|default('') }}"
    foreach ($p in $d.Paragraphs) {
        $ptext = $p.Range.Text.TrimEnd([char]13)
        if ($ptext -eq $needle) {
            $rng = $p.Range
            $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">&#8226; </w:t></w:r><w:r><w:t>&#8226;</w:t></w:r></w:p>'
            $rng.InsertXML($xml)
            break
        }
    }
}

# Second pass: give each of the now-updated bullet paragraphs the explicit
# indentation (w:ind w:left="360" w:firstLine="0") that the diff calls for.
# This must be done via ParagraphFormat (rather than embedded in the
# InsertXML payload) because InsertXML does not persist <w:ind>.
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "$([char]0x2022) $([char]0x2022)") {
        $p.Range.ParagraphFormat.LeftIndent = 18
        $p.Range.ParagraphFormat.FirstLineIndent = 0
    }
}
